$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.054.80"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "1.689.21"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.14%  "
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0627"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.927.80"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "1.689.01"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.19"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.558"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "251.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.59%  "
$ws.Range("D18").Value = "28.009.40"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.24"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.49"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +7.02%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "1.431.88"
$ws.Range("E34").Value = "  -7.07%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.942"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.592"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("D44").Value = "1.835.12"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.24"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("E47").Value = "  +6.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "0.0₆0112"
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.85"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.65%  "
